# Apply the "Updated cryptos list" refresh to the crypto tracker sheet.
# Values that look like plain numbers (e.g. "385.06") must be forced back
# to Text so Excel keeps them as strings instead of silently converting
# them to numeric cells (the source data uses Text-formatted price cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '51.586.97'
$ws.Range('E2').Value = '  +1.10%  '

# Row 3
$ws.Range('D3').Value = '3.034.34'
$ws.Range('E3').Value = '  +2.69%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '385.06'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.43%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.36'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.50%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.587'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.35%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.73'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.31%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.137'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.01%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0859'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.48%  '

# Row 13
$ws.Range('D13').Value = '3.514.28'
$ws.Range('E13').Value = '  +2.82%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.65'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.93%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.74'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.22%  '

# Row 16
$ws.Range('D16').Value = '3.037.75'
$ws.Range('E16').Value = '  +2.65%  '

# Row 17
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.974'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.72%  '

# Row 18
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.67'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -11.04%  '

# Row 19
$ws.Range('D19').Value = '51.614.66'
$ws.Range('E19').Value = '  +1.03%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.09'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.19%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.41'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.35%  '

# Row 22
$ws.Range('E22').Value = '  +0.18%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.31%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.05'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.19%  '

# Row 25
$ws.Range('E25').Value = '  -3.53%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.31'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +5.34%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.47'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.79%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.173'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.45%  '

# Row 29
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.10%  '

# Row 30
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.24'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.61%  '

# Row 31
$ws.Range('E31').Value = '  -1.19%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.28'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.10%  '

# Row 33
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.53'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.91%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '33.88'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.08%  '

# Row 35
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.05'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.70%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0448'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.87%  '

# Row 37
$ws.Range('E37').Value = '  -0.14%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.32'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.44%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.289'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.94%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.98'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.55%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.86'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.58%  '

# Row 42
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '128.47'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.19%  '

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.116'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.05%  '

# Row 44
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.52'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.24%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.68'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.35%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.61'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.33%  '

# Row 47
$ws.Range('E47').Value = '  +2.49%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.43'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.96%  '

# Row 49
$ws.Range('D49').Value = '2.026.87'
$ws.Range('E49').Value = '  -1.28%  '

# Row 50
$ws.Range('D50').Value = '3.334.53'
$ws.Range('E50').Value = '  +2.58%  '

# Row 51
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0317'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.51%  '
